$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "2022-08-16"
$ws.Range("J2").Value = 90
$ws.Range("K2").Value = 24000
$ws.Range("L2").Value = 24000
$ws.Range("M2").Value = 24000
$ws.Range("P2").Value = 1600
$ws.Range("D3").Value = "2022-10-04"
$ws.Range("J3").Value = 80
$ws.Range("K3").Value = 22000
$ws.Range("L3").Value = 22000
$ws.Range("M3").Value = 22000
$ws.Range("P3").Value = 1467
$ws.Range("D4").Value = "2022-07-22"
$ws.Range("D6").Value = "2021-06-18"
$ws.Range("J6").Value = 80
$ws.Range("K6").Value = 25000
$ws.Range("L6").Value = 25000
$ws.Range("M6").Value = 25000
$ws.Range("P6").Value = 1667
$ws.Range("D7").Value = "2022-08-26"
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 23000
$ws.Range("L7").Value = 23000
$ws.Range("M7").Value = 23000
$ws.Range("P7").Value = 1533
$ws.Range("D8").Value = "2022-07-12"
$ws.Range("J8").Value = 90
$ws.Range("K8").Value = 25000
$ws.Range("L8").Value = 25000
$ws.Range("M8").Value = 25000
$ws.Range("P8").Value = 1667
$ws.Range("D9").Value = "2022-07-08"
$ws.Range("K9").Value = 25000
$ws.Range("L9").Value = 25000
$ws.Range("M9").Value = 25000
$ws.Range("P9").Value = 1667
$ws.Range("D10").Value = "2022-07-15"
$ws.Range("D11").Value = "2022-09-27"
$ws.Range("D12").Value = "2022-08-05"
$ws.Range("J12").Value = 120
$ws.Range("K12").Value = 24000
$ws.Range("L12").Value = 24000
$ws.Range("M12").Value = 24000
$ws.Range("P12").Value = 1600
$ws.Range("D13").Value = "2022-08-02"
$ws.Range("D14").Value = "2022-09-23"
$ws.Range("J14").Value = 90
$ws.Range("D15").Value = "2022-09-13"
$ws.Range("K15").Value = 23000
$ws.Range("L15").Value = 23000
$ws.Range("M15").Value = 23000
$ws.Range("P15").Value = 1533
$ws.Range("D16").Value = "2021-07-23"
$ws.Range("J16").Value = 80
$ws.Range("K16").Value = 25000
$ws.Range("L16").Value = 25000
$ws.Range("M16").Value = 25000
$ws.Range("P16").Value = 1667
$ws.Range("D17").Value = "2022-07-19"
$ws.Range("J17").Value = 100
$ws.Range("L17").Value = 25000
$ws.Range("M17").Value = 24000
$ws.Range("P17").Value = 1600
$ws.Range("D18").Value = "2022-08-09"
$ws.Range("J18").Value = 120
$ws.Range("K18").Value = 24000
$ws.Range("L18").Value = 24000
$ws.Range("M18").Value = 24000
$ws.Range("P18").Value = 1600
$ws.Range("D19").Value = "2022-07-29"
$ws.Range("K19").Value = 25000
$ws.Range("L19").Value = 25000
$ws.Range("M19").Value = 25000
$ws.Range("P19").Value = 1667
$ws.Range("D20").Value = "2022-08-19"
$ws.Range("J20").Value = 120
$ws.Range("K20").Value = 24000
$ws.Range("L20").Value = 24000
$ws.Range("M20").Value = 24000
$ws.Range("P20").Value = 1600
$ws.Range("D21").Value = "2021-08-10"
$ws.Range("J21").Value = 90
$ws.Range("K21").Value = 25000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 25000
$ws.Range("P21").Value = 1667
$ws.Range("D22").Value = "2022-08-30"
$ws.Range("K22").Value = 24000
$ws.Range("L22").Value = 24000
$ws.Range("M22").Value = 24000
$ws.Range("P22").Value = 1600
$ws.Range("D23").Value = "2022-08-08"
$ws.Range("J23").Value = 70
$ws.Range("D24").Value = "2022-06-28"
$ws.Range("D25").Value = "2021-07-30"
$ws.Range("J25").Value = 90
$ws.Range("K25").Value = 25000
$ws.Range("L25").Value = 25000
$ws.Range("M25").Value = 25000
$ws.Range("P25").Value = 1667
$ws.Range("D26").Value = "2022-09-15"
$ws.Range("J26").Value = 70
$ws.Range("K26").Value = 22000
$ws.Range("L26").Value = 22000
$ws.Range("M26").Value = 22000
$ws.Range("P26").Value = 1467
$ws.Range("D27").Value = "2022-09-02"
$ws.Range("J27").Value = 70
$ws.Range("K27").Value = 23000
$ws.Range("L27").Value = 23000
$ws.Range("M27").Value = 23000
$ws.Range("P27").Value = 1533
